# Apply "added alphapose and qol" changes to Results.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
# A1 "Unnamed: 0" and B1 "typ filmu" stay the same.
$ws.Range("C1").Value = "Brakujące klatki"
$ws.Range("D1").Value = "Poprawność klasyfikatora"
$ws.Range("E1").Value = "ilość punktów"

# --- Clear out the old data body (rows 2:19, cols A:E) ---
$ws.Range("A2:E19").ClearContents()

# Column E no longer holds any data under the new layout - only the header remains.

# --- Write the new data body ---
# Each element: model label (col A, only on first sub-row), metric label (col B),
# "Brakujące klatki" count (col C), classifier correctness score (col D)
$data = @(
    @("BlazePose", "f1 norma",    9,   0.8910318225650916),
    @("",          "f1 ciemno",   3,   0.6730954676952748),
    @("",          "f1 zakrycie", 244, 0.5332690453230472),
    @("",          "f1 oba",      367, 0.4503375120540019),
    @("",          "f2 norma",    1,   0.7129543336439889),
    @("",          "f2 ciemno",   18,  0.5983224603914259),
    @("",          "f2 zakrycie", 131, 0.6132339235787512),
    @("",          "f2 oba",      298, 0.5955265610438024),
    @("OpenPose",  "f1 norma",    13,  0.8119575699132112),
    @("",          "f1 ciemno",   2,   0.5805207328833173),
    @("",          "f1 zakrycie", 58,  0.5940212150433944),
    @("",          "f1 oba",      101, 0.5332690453230472),
    @("",          "f2 norma",    0,   0.7502329916123019),
    @("",          "f2 ciemno",   15,  0.5517241379310345),
    @("",          "f2 zakrycie", 2,   0.527493010251631),
    @("",          "f2 oba",      90,  0.5405405405405406),
    @("AlphaPose", "f1 norma",    298, 0.8871745419479267),
    @("",          "f1 ciemno",   1,   0.5226615236258437),
    @("",          "f1 zakrycie", 1,   0.3278688524590164),
    @("",          "f1 oba",      1,   0.5429122468659595),
    @("",          "f2 norma",    1,   0.6728797763280522),
    @("",          "f2 ciemno",   1,   0.4986020503261883),
    @("",          "f2 zakrycie", 1,   0.4986020503261883),
    @("",          "f2 oba",      1,   0.4958061509785648)
)

$row = 2
foreach ($entry in $data) {
    if ($entry[0] -ne "") {
        $ws.Cells.Item($row, 1).Value = $entry[0]
    }
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
